# Add new hire rows 126-156 to the qabul worksheet (F.I.SH, Ta'lim yo'nalishi, Passport,
# Shartnoma raqam, Viloyat, Tuman, Telefon raqam, Sana) -- data added 03-12-2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(126, 1).Value = "Ramatov Yangiboy Jumanazarovich"
$ws.Cells.Item(126, 2).Value = "Maktabgacha talim tashkiloti musiqa rahbari"
$ws.Cells.Item(126, 3).Value = "AB2864318"
$ws.Cells.Item(126, 4).Value = "'301"
$ws.Cells.Item(126, 5).Value = "Xorazm viloyati"
$ws.Cells.Item(126, 6).Value = "Urganch tumani"
$ws.Cells.Item(126, 7).Value = "'998912784165"
$ws.Cells.Item(126, 8).Value = "'28-11-2024"

$ws.Cells.Item(127, 1).Value = "Sadilloyeva Dilfuza Shuxratovna"
$ws.Cells.Item(127, 2).Value = "Maktabgacha talim tashkiloti defektolog/logopedi 576 soat"
$ws.Cells.Item(127, 3).Value = "AB1130811"
$ws.Cells.Item(127, 4).Value = "'302"
$ws.Cells.Item(127, 5).Value = "Buxoro viloyati"
$ws.Cells.Item(127, 6).Value = "Kogon tumani"
$ws.Cells.Item(127, 7).Value = "'998995558780"
$ws.Cells.Item(127, 8).Value = "'28-11-2024"

$ws.Cells.Item(128, 1).Value = "Babayeva Kamola Xusanovna"
$ws.Cells.Item(128, 2).Value = "Maktabgacha talim tashkiloti defektolog/logopedi 576 soat"
$ws.Cells.Item(128, 3).Value = "AB4404861"
$ws.Cells.Item(128, 4).Value = "'303"
$ws.Cells.Item(128, 5).Value = "Toshkent viloyati"
$ws.Cells.Item(128, 6).Value = "Angren shahri"
$ws.Cells.Item(128, 7).Value = "'998931852007"
$ws.Cells.Item(128, 8).Value = "'28-11-2024"

$ws.Cells.Item(129, 1).Value = "Sadilloyeva Dilfuza Shuxratovna"
$ws.Cells.Item(129, 2).Value = "Maktabgacha talim tashkiloti defektolog/logopedi"
$ws.Cells.Item(129, 3).Value = "AB1130811"
$ws.Cells.Item(129, 4).Value = "'304"
$ws.Cells.Item(129, 5).Value = "Buxoro viloyati"
$ws.Cells.Item(129, 6).Value = "Buxoro tumani"
$ws.Cells.Item(129, 7).Value = "'998995558780"
$ws.Cells.Item(129, 8).Value = "'28-11-2024"

$ws.Cells.Item(130, 1).Value = "Jurayeva Raxima Habibullayevna"
$ws.Cells.Item(130, 2).Value = "Maktabgacha talim tashkiloti direktori"
$ws.Cells.Item(130, 3).Value = "AD3610244"
$ws.Cells.Item(130, 4).Value = "'305"
$ws.Cells.Item(130, 5).Value = "Namangan viloyati"
$ws.Cells.Item(130, 6).Value = "Kosonsoy tumani"
$ws.Cells.Item(130, 7).Value = "'998951015074"
$ws.Cells.Item(130, 8).Value = "'28-11-2024"

$ws.Cells.Item(131, 1).Value = "Yo'ldasheva Nilufar Ibrohimovna"
$ws.Cells.Item(131, 2).Value = "Maktabgacha talim tashkiloti direktori"
$ws.Cells.Item(131, 3).Value = "AD8567923"
$ws.Cells.Item(131, 4).Value = "'306"
$ws.Cells.Item(131, 5).Value = "Andijon viloyati"
$ws.Cells.Item(131, 6).Value = "Andijon shahri"
$ws.Cells.Item(131, 7).Value = "'998916004616"
$ws.Cells.Item(131, 8).Value = "'28-11-2024"

$ws.Cells.Item(132, 1).Value = "Botirova Zamira O'ralovna"
$ws.Cells.Item(132, 2).Value = "Maktabgacha talim tashkiloti tarbiyachisi"
$ws.Cells.Item(132, 3).Value = "AD8925197"
$ws.Cells.Item(132, 4).Value = "'307"
$ws.Cells.Item(132, 5).Value = "Toshkent shahri"
$ws.Cells.Item(132, 6).Value = "Mirzo Ulugʻbek tumani"
$ws.Cells.Item(132, 7).Value = "'998909281743"
$ws.Cells.Item(132, 8).Value = "'28-11-2024"

$ws.Cells.Item(133, 1).Value = "Ravshanova Marjona Aminjonovna"
$ws.Cells.Item(133, 2).Value = "Maktabgacha talim tashkiloti tarbiyachisi 576 soat"
$ws.Cells.Item(133, 3).Value = "AB4649674"
$ws.Cells.Item(133, 4).Value = "'308"
$ws.Cells.Item(133, 5).Value = "Samarqand viloyati"
$ws.Cells.Item(133, 6).Value = "Samarqand tumani"
$ws.Cells.Item(133, 7).Value = "'998906564240"
$ws.Cells.Item(133, 8).Value = "'29-11-2024"

$ws.Cells.Item(134, 1).Value = "Guylmamedova Raisa Viktorovna"
$ws.Cells.Item(134, 2).Value = "Maktabgacha talim tashkiloti tarbiyachisi 576 soat"
$ws.Cells.Item(134, 3).Value = "AD6571021"
$ws.Cells.Item(134, 4).Value = "'309"
$ws.Cells.Item(134, 5).Value = "Toshkent shahri"
$ws.Cells.Item(134, 6).Value = "Chilonzor tumani"
$ws.Cells.Item(134, 7).Value = "'998994805061"
$ws.Cells.Item(134, 8).Value = "'29-11-2024"

$ws.Cells.Item(135, 1).Value = "Komilova Nodiraxon Mahmudjon qizi"
$ws.Cells.Item(135, 2).Value = "Maktabgacha ta``lim tashkiloti tashkilot oshpazi"
$ws.Cells.Item(135, 3).Value = "AB5582671"
$ws.Cells.Item(135, 4).Value = "'310"
$ws.Cells.Item(135, 5).Value = "Andijon viloyati"
$ws.Cells.Item(135, 6).Value = "Paxtaobod tumani"
$ws.Cells.Item(135, 7).Value = "'998948893272"
$ws.Cells.Item(135, 8).Value = "'29-11-2024"

$ws.Cells.Item(136, 1).Value = "Mirzamatova Zulxumor Umaraliyevna"
$ws.Cells.Item(136, 2).Value = "Maktabgacha ta``lim tashkiloti tashkilot oshpazi"
$ws.Cells.Item(136, 3).Value = "AD4202739"
$ws.Cells.Item(136, 4).Value = "'311"
$ws.Cells.Item(136, 5).Value = "Andijon viloyati"
$ws.Cells.Item(136, 6).Value = "Paxtaobod tumani"
$ws.Cells.Item(136, 7).Value = "'998936322172"
$ws.Cells.Item(136, 8).Value = "'29-11-2024"

$ws.Cells.Item(137, 1).Value = "Isakova Ozodaxon Tulashovna"
$ws.Cells.Item(137, 2).Value = "Maktabgacha talim tashkiloti tarbiyachisi"
$ws.Cells.Item(137, 3).Value = "AA9447846"
$ws.Cells.Item(137, 4).Value = "'312"
$ws.Cells.Item(137, 5).Value = "Toshkent shahri"
$ws.Cells.Item(137, 6).Value = "Mirzo Ulugʻbek tumani"
$ws.Cells.Item(137, 7).Value = "'998950150201"
$ws.Cells.Item(137, 8).Value = "'29-11-2024"

$ws.Cells.Item(138, 1).Value = "Maxkamova Ra'no Sadikjanovna"
$ws.Cells.Item(138, 2).Value = "Maktabgacha talim tashkiloti direktori"
$ws.Cells.Item(138, 3).Value = "AD7811878"
$ws.Cells.Item(138, 4).Value = "'313"
$ws.Cells.Item(138, 5).Value = "Toshkent viloyati"
$ws.Cells.Item(138, 6).Value = "Yuqori Chirchiq tumani"
$ws.Cells.Item(138, 7).Value = "'+998944298109"
$ws.Cells.Item(138, 8).Value = "'30-11-2024"

$ws.Cells.Item(139, 1).Value = "Usmanova Zulfiya Rustamovna"
$ws.Cells.Item(139, 2).Value = "Maktabgacha talim tashkiloti direktori"
$ws.Cells.Item(139, 3).Value = "AD1510147"
$ws.Cells.Item(139, 4).Value = "'314"
$ws.Cells.Item(139, 5).Value = "Fargona viloyati"
$ws.Cells.Item(139, 6).Value = "Fargʻona tumani"
$ws.Cells.Item(139, 7).Value = "'998948364949"
$ws.Cells.Item(139, 8).Value = "'30-11-2024"

$ws.Cells.Item(140, 1).Value = "Mutalipova Nargiza Zakirovna"
$ws.Cells.Item(140, 2).Value = "Maktabgacha talim tashkiloti defektolog/logopedi"
$ws.Cells.Item(140, 3).Value = "AD8099504"
$ws.Cells.Item(140, 4).Value = "'315"
$ws.Cells.Item(140, 5).Value = "Toshkent shahri"
$ws.Cells.Item(140, 6).Value = "Yunusobod tumani"
$ws.Cells.Item(140, 7).Value = "'+998944298109"
$ws.Cells.Item(140, 8).Value = "'30-11-2024"

$ws.Cells.Item(141, 1).Value = "Mahmudova Aziza Mirzahamdamovna"
$ws.Cells.Item(141, 2).Value = "Maktabgacha talim tashkiloti tarbiyachisi"
$ws.Cells.Item(141, 3).Value = "AD7813109"
$ws.Cells.Item(141, 4).Value = "'316"
$ws.Cells.Item(141, 5).Value = "Toshkent shahri"
$ws.Cells.Item(141, 6).Value = "Yunusobod tumani"
$ws.Cells.Item(141, 7).Value = "'998777011980"
$ws.Cells.Item(141, 8).Value = "'30-11-2024"

$ws.Cells.Item(142, 1).Value = "Avazova Muxayyoxon Muzaffarjon qizi"
$ws.Cells.Item(142, 2).Value = "Maktabgacha talim tashkiloti metodisti"
$ws.Cells.Item(142, 3).Value = "AC1436252"
$ws.Cells.Item(142, 4).Value = "'317"
$ws.Cells.Item(142, 5).Value = "Namangan viloyati"
$ws.Cells.Item(142, 6).Value = "Pop tumani"
$ws.Cells.Item(142, 7).Value = "'998976228333"
$ws.Cells.Item(142, 8).Value = "'30-11-2024"

$ws.Cells.Item(143, 1).Value = "Xayrullayeva Shohista Bekmirza qizi"
$ws.Cells.Item(143, 2).Value = "Maktabgacha talim tashkiloti metodisti"
$ws.Cells.Item(143, 3).Value = "AC0459728"
$ws.Cells.Item(143, 4).Value = "'318"
$ws.Cells.Item(143, 5).Value = "Namangan viloyati"
$ws.Cells.Item(143, 6).Value = "Yangi Namangan"
$ws.Cells.Item(143, 7).Value = "'998949709559"
$ws.Cells.Item(143, 8).Value = "'01-12-2024"

$ws.Cells.Item(144, 1).Value = "Nurullayeva Nodira Maxmudovna"
$ws.Cells.Item(144, 2).Value = "Maktabgacha talim tashkiloti direktori"
$ws.Cells.Item(144, 3).Value = "AA8228491"
$ws.Cells.Item(144, 4).Value = "'319"
$ws.Cells.Item(144, 5).Value = "Buxoro viloyati"
$ws.Cells.Item(144, 6).Value = "Buxoro tumani"
$ws.Cells.Item(144, 7).Value = "'998997361179"
$ws.Cells.Item(144, 8).Value = "'02-12-2024"

$ws.Cells.Item(145, 1).Value = "YULDASHEVA DILDORAHON"
$ws.Cells.Item(145, 2).Value = "Maktabgacha talim tashkiloti direktori"
$ws.Cells.Item(145, 3).Value = "AA9049494"
$ws.Cells.Item(145, 4).Value = "'320"
$ws.Cells.Item(145, 5).Value = "Andijon viloyati"
$ws.Cells.Item(145, 6).Value = "Andijon shahri"
$ws.Cells.Item(145, 7).Value = "'998932594176"
$ws.Cells.Item(145, 8).Value = "'02-12-2024"

$ws.Cells.Item(146, 1).Value = "Irgasheva Shalola Djabbarovna"
$ws.Cells.Item(146, 2).Value = "Maktabgacha talim tashkiloti direktori"
$ws.Cells.Item(146, 3).Value = "AD1330938"
$ws.Cells.Item(146, 4).Value = "'321"
$ws.Cells.Item(146, 5).Value = "Samarqand viloyati"
$ws.Cells.Item(146, 6).Value = "Samarqand shahri"
$ws.Cells.Item(146, 7).Value = "'998992912547"
$ws.Cells.Item(146, 8).Value = "'02-12-2024"

$ws.Cells.Item(147, 1).Value = "Islomova Mahliyo Komil qizi"
$ws.Cells.Item(147, 2).Value = "Maktabgacha talim tashkiloti psixologi"
$ws.Cells.Item(147, 3).Value = "AC0282022"
$ws.Cells.Item(147, 4).Value = "'322"
$ws.Cells.Item(147, 5).Value = "Jizzax viloyati"
$ws.Cells.Item(147, 6).Value = "Sharof Rashidov tumani"
$ws.Cells.Item(147, 7).Value = "'998932918995"
$ws.Cells.Item(147, 8).Value = "'02-12-2024"

$ws.Cells.Item(148, 1).Value = "Rasulova Nargiza Raximdjanovna"
$ws.Cells.Item(148, 2).Value = "Maktabgacha talim tashkiloti direktori"
$ws.Cells.Item(148, 3).Value = "AD6566738"
$ws.Cells.Item(148, 4).Value = "'323"
$ws.Cells.Item(148, 5).Value = "Toshkent shahri"
$ws.Cells.Item(148, 6).Value = "Chilonzor tumani"
$ws.Cells.Item(148, 7).Value = "'998909313092"
$ws.Cells.Item(148, 8).Value = "'02-12-2024"

$ws.Cells.Item(149, 1).Value = "Bozorova Aziz Djuraqulovna"
$ws.Cells.Item(149, 2).Value = "Maktabgacha talim tashkiloti direktori"
$ws.Cells.Item(149, 3).Value = "AB5048108"
$ws.Cells.Item(149, 4).Value = "'324"
$ws.Cells.Item(149, 5).Value = "Samarqand viloyati"
$ws.Cells.Item(149, 6).Value = "Samarqand shahri"
$ws.Cells.Item(149, 7).Value = "'998915552375"
$ws.Cells.Item(149, 8).Value = "'02-12-2024"

$ws.Cells.Item(150, 1).Value = "Ahmedova Muhabbatxon Xojiakbarxon qizi"
$ws.Cells.Item(150, 2).Value = "Maktabgacha talim tashkiloti psixologi"
$ws.Cells.Item(150, 3).Value = "AD6646697"
$ws.Cells.Item(150, 4).Value = "'325"
$ws.Cells.Item(150, 5).Value = "Toshkent shahri"
$ws.Cells.Item(150, 6).Value = "Shayxontohur tumani"
$ws.Cells.Item(150, 7).Value = "'998950775055"
$ws.Cells.Item(150, 8).Value = "'02-12-2024"

$ws.Cells.Item(151, 1).Value = "Urunova Dildora Sirojiddinovna"
$ws.Cells.Item(151, 2).Value = "Maktabgacha talim tashkiloti direktori"
$ws.Cells.Item(151, 3).Value = "AD1866387"
$ws.Cells.Item(151, 4).Value = "'326"
$ws.Cells.Item(151, 5).Value = "Samarqand viloyati"
$ws.Cells.Item(151, 6).Value = "Samarqand shahri"
$ws.Cells.Item(151, 7).Value = "'998982732145"
$ws.Cells.Item(151, 8).Value = "'02-12-2024"

$ws.Cells.Item(152, 1).Value = "Xudoyberdiyeva Zarnigor Kobiljonovna"
$ws.Cells.Item(152, 2).Value = "Maktabgacha talim tashkiloti psixologi"
$ws.Cells.Item(152, 3).Value = "AA9000453"
$ws.Cells.Item(152, 4).Value = "'327"
$ws.Cells.Item(152, 5).Value = "Samarqand viloyati"
$ws.Cells.Item(152, 6).Value = "Samarqand tumani"
$ws.Cells.Item(152, 7).Value = "'998979271393"
$ws.Cells.Item(152, 8).Value = "'02-12-2024"

$ws.Cells.Item(153, 1).Value = "Muhammadieva Shoira Holmuratovna"
$ws.Cells.Item(153, 2).Value = "Maktabgacha talim tashkiloti direktori"
$ws.Cells.Item(153, 3).Value = "AD4160337"
$ws.Cells.Item(153, 4).Value = "'328"
$ws.Cells.Item(153, 5).Value = "Samarqand viloyati"
$ws.Cells.Item(153, 6).Value = "Samarqand shahri"
$ws.Cells.Item(153, 7).Value = "'998933552525"
$ws.Cells.Item(153, 8).Value = "'02-12-2024"

$ws.Cells.Item(154, 1).Value = "Ochilova Dildora Normamatovna"
$ws.Cells.Item(154, 2).Value = "Maktabgacha talim tashkiloti psixologi"
$ws.Cells.Item(154, 3).Value = "AD4285763"
$ws.Cells.Item(154, 4).Value = "'329"
$ws.Cells.Item(154, 5).Value = "Qashqadaryo viloyati"
$ws.Cells.Item(154, 6).Value = "Qarshi tumani"
$ws.Cells.Item(154, 7).Value = "'998973146151"
$ws.Cells.Item(154, 8).Value = "'02-12-2024"

$ws.Cells.Item(155, 1).Value = "Qodirova Dilorom Alivayevna"
$ws.Cells.Item(155, 2).Value = "Maktabgacha talim tashkiloti tarbiyachisi"
$ws.Cells.Item(155, 3).Value = "AA2991834"
$ws.Cells.Item(155, 4).Value = "'330"
$ws.Cells.Item(155, 5).Value = "Toshkent shahri"
$ws.Cells.Item(155, 6).Value = "Yunusobod tumani"
$ws.Cells.Item(155, 7).Value = "'998908058073"
$ws.Cells.Item(155, 8).Value = "'02-12-2024"

$ws.Cells.Item(156, 1).Value = "Karimova Bonu Nazirjon qizi"
$ws.Cells.Item(156, 2).Value = "Maktabgacha talim tashkiloti tarbiyachisi"
$ws.Cells.Item(156, 3).Value = "AB3171517"
$ws.Cells.Item(156, 4).Value = "'331"
$ws.Cells.Item(156, 5).Value = "Toshkent shahri"
$ws.Cells.Item(156, 6).Value = "Yunusobod tumani"
$ws.Cells.Item(156, 7).Value = "'998958154892"
$ws.Cells.Item(156, 8).Value = "'03-12-2024"

